# Refresh the cryptos list (GitHub Actions scheduled update).
# Updates Price (D) / Volume(1h) (E) text values for most rows, and
# fully replaces row 51 (dogwifhat -> SuiNetwork).
#
# NOTE: columns D/E store plain-looking numbers/percentages as TEXT
# (inline strings) in the source workbook. Excel's Range.Value setter
# auto-coerces numeric-looking strings ("582.54") into real numbers
# when the cell's NumberFormat is "General". To keep those cells text
# (matching the original data type), we force NumberFormat = "@" just
# before writing any value that would otherwise parse as a number.
# Values that contain extra punctuation (two dots, e.g. "64.328.30"),
# or that include the surrounding spaces / a trailing "%", already fail
# Excel's numeric parser, so no format change is needed for those.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "64.328.30"
$ws.Range("E2").Value = "  +5.66%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.476.81"
$ws.Range("E3").Value = "  +7.07%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.10%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.54"
$ws.Range("E5").Value = "  +6.81%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.59"
$ws.Range("E6").Value = "  +7.61%  "

# Row 7 - USDC
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.09%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.482.25"
$ws.Range("E8").Value = "  +7.07%  "

# Row 9 - XRP
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.536"
$ws.Range("E9").Value = "  +1.37%  "

# Row 10 - Toncoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.61"
$ws.Range("E10").Value = "  +2.98%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +7.80%  "

# Row 12 - Cardano
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.442"
$ws.Range("E12").Value = "  +1.97%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "4.083.12"
$ws.Range("E13").Value = "  +7.23%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -0.27%  "

# Row 15 - ShibaInu
$ws.Range("E15").Value = "  +7.82%  "

# Row 16 - Avalanche
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.73"
$ws.Range("E16").Value = "  +4.88%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "64.502.73"
$ws.Range("E17").Value = "  +6.03%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.488.52"
$ws.Range("E18").Value = "  +6.38%  "

# Row 19 - Polkadot
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.45"
$ws.Range("E19").Value = "  +2.06%  "

# Row 20 - Chainlink
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.40"
$ws.Range("E20").Value = "  +7.02%  "

# Row 21 - BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "398.96"
$ws.Range("E21").Value = "  +5.64%  "

# Row 22 - Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.54"
$ws.Range("E22").Value = "  +1.23%  "

# Row 23 - Polygon
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.547"
$ws.Range("E23").Value = "  +2.53%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  +0.04%  "

# Row 25 - Litecoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.09"
$ws.Range("E25").Value = "  +2.96%  "

# Row 26 - PEPE
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000111"
$ws.Range("E26").Value = "  +20.51%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.45"
$ws.Range("E27").Value = "  +9.69%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  +6.58%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("E29").Value = "  -0.26%  "

# Row 30 - Fetch.AI
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.41"
$ws.Range("E30").Value = "  +13.94%  "

# Row 31 - NEARProtocol
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.96"
$ws.Range("E31").Value = "  +10.03%  "

# Row 32 - RenderToken
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.73"
$ws.Range("E32").Value = "  +8.29%  "

# Row 33 - PancakeSwap
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.05"
$ws.Range("E33").Value = "  +6.17%  "

# Row 34 - EthereumClassic
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.89"
$ws.Range("E34").Value = "  +5.64%  "

# Row 36 - Aptos
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.95"
$ws.Range("E36").Value = "  +4.45%  "

# Row 37 - ImmutableX
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.51"
$ws.Range("E37").Value = "  +5.82%  "

# Row 38 - Monero
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "158.71"
$ws.Range("E38").Value = "  -0.14%  "

# Row 39 - EnergySwap
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "28.55"
$ws.Range("E39").Value = "  +7.82%  "

# Row 40 - Hedera
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0790"
$ws.Range("E40").Value = "  +9.17%  "

# Row 41 - Stacks
$ws.Range("E41").Value = "  +9.57%  "

# Row 42 - Maker
$ws.Range("D42").Value = "2.894.68"
$ws.Range("E42").Value = "  +3.30%  "

# Row 43 - VeChain
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0324"
$ws.Range("E43").Value = "  +3.15%  "

# Row 44 - Mantle
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.784"
$ws.Range("E44").Value = "  +7.30%  "

# Row 45 - Filecoin
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.44"
$ws.Range("E45").Value = "  +3.41%  "

# Row 46 - OKB
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "42.15"
$ws.Range("E46").Value = "  +5.48%  "

# Row 47 - ONDO
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.11"
$ws.Range("E47").Value = "  +10.43%  "

# Row 48 - RenzoRestakedETH
$ws.Range("D48").Value = "3.532.05"
$ws.Range("E48").Value = "  +7.40%  "

# Row 49 - InjectiveProtocol
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.81"
$ws.Range("E49").Value = "  +5.79%  "

# Row 50 - Bittensor
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "299.85"
$ws.Range("E50").Value = "  +8.01%  "

# Row 51 - dogwifhat -> SuiNetwork (full row replacement)
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.841"
$ws.Range("E51").Value = "  +4.46%  "
